$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Range("D11").Value = "TIMESTAMP"
$ws.Range("D13").Value = "TIMESTAMP"
[void]$ws.Range("D11").Select()
